# Update "paises.xlsx" (Pais sheet) with refreshed COVID-19 country stats
# and corrected country ordering/labels, per upstream data refresh
# (commit: "Update countries & provincias Spain").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (A1): "11:28" -> "12:45"
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 8 de Julio de 2020 a las 12:45"

# Row 4: refreshed totals
$ws.Cells.Item(4, 2).Value = 3097538
$ws.Cells.Item(4, 3).Value = 454
$ws.Cells.Item(4, 4).Value = 1355524
$ws.Cells.Item(4, 5).Value = 1608023
$ws.Cells.Item(4, 7).Value = 19
$ws.Cells.Item(4, 8).Value = 133991

# Row 6: refreshed totals
$ws.Cells.Item(6, 2).Value = 746500
$ws.Cells.Item(6, 3).Value = 3019
$ws.Cells.Item(6, 4).Value = 458615
$ws.Cells.Item(6, 5).Value = 267201
$ws.Cells.Item(6, 7).Value = 31
$ws.Cells.Item(6, 8).Value = 20684

# Row 13: refreshed totals
$ws.Cells.Item(13, 2).Value = 248379
$ws.Cells.Item(13, 3).Value = 2691
$ws.Cells.Item(13, 4).Value = 209463
$ws.Cells.Item(13, 5).Value = 26832
$ws.Cells.Item(13, 7).Value = 153
$ws.Cells.Item(13, 8).Value = 12084

# Row 35: refreshed totals
$ws.Cells.Item(35, 2).Value = 52007
$ws.Cells.Item(35, 3).Value = 762
$ws.Cells.Item(35, 4).Value = 42108
$ws.Cells.Item(35, 5).Value = 9520
$ws.Cells.Item(35, 7).Value = 2
$ws.Cells.Item(35, 8).Value = 379

# Row 51: now "Rumania"
$ws.Cells.Item(51, 1).Value = "Rumania"
$ws.Cells.Item(51, 2).Value = 30175
$ws.Cells.Item(51, 3).Value = 555
$ws.Cells.Item(51, 4).Value = 20799
$ws.Cells.Item(51, 5).Value = 7559
$ws.Cells.Item(51, 8).Value = 1817

# Row 52: now "Armenia"
$ws.Cells.Item(52, 1).Value = "Armenia"
$ws.Cells.Item(52, 2).Value = 29820
$ws.Cells.Item(52, 3).Value = 535
$ws.Cells.Item(52, 4).Value = 17427
$ws.Cells.Item(52, 5).Value = 11872
$ws.Cells.Item(52, 7).Value = 18
$ws.Cells.Item(52, 8).Value = 521

# Row 53: now "Nigeria"
$ws.Cells.Item(53, 1).Value = "Nigeria"
$ws.Cells.Item(53, 2).Value = 29789
$ws.Cells.Item(53, 4).Value = 12108
$ws.Cells.Item(53, 5).Value = 17012
$ws.Cells.Item(53, 8).Value = 669

# Row 60: refreshed totals
$ws.Cells.Item(60, 2).Value = 18513
$ws.Cells.Item(60, 3).Value = 92
$ws.Cells.Item(60, 4).Value = 16721
$ws.Cells.Item(60, 5).Value = 1086

# Row 75: refreshed totals
$ws.Cells.Item(75, 2).Value = 8677
$ws.Cells.Item(75, 3).Value = 3
$ws.Cells.Item(75, 4).Value = 8486
$ws.Cells.Item(75, 5).Value = 70

# Row 77: now "Kenia"
$ws.Cells.Item(77, 1).Value = "Kenia"
$ws.Cells.Item(77, 2).Value = 8528
$ws.Cells.Item(77, 3).Value = 278
$ws.Cells.Item(77, 4).Value = 2593
$ws.Cells.Item(77, 5).Value = 5766
$ws.Cells.Item(77, 7).Value = 2
$ws.Cells.Item(77, 8).Value = 169

# Row 78: now "Kirguistan"
$ws.Cells.Item(78, 1).Value = "Kirguistan"
$ws.Cells.Item(78, 2).Value = 8486
$ws.Cells.Item(78, 3).Value = 345
$ws.Cells.Item(78, 4).Value = 2983
$ws.Cells.Item(78, 5).Value = 5391
$ws.Cells.Item(78, 7).Value = 13
$ws.Cells.Item(78, 8).Value = 112

# Row 80: refreshed totals
$ws.Cells.Item(80, 2).Value = 7657
$ws.Cells.Item(80, 3).Value = 110
$ws.Cells.Item(80, 4).Value = 5097
$ws.Cells.Item(80, 5).Value = 2419
$ws.Cells.Item(80, 7).Value = 4
$ws.Cells.Item(80, 8).Value = 141

# Row 103: now "Albania"
$ws.Cells.Item(103, 1).Value = "Albania"
$ws.Cells.Item(103, 2).Value = 3106
$ws.Cells.Item(103, 3).Value = 68
$ws.Cells.Item(103, 4).Value = 1791
$ws.Cells.Item(103, 5).Value = 1232
$ws.Cells.Item(103, 7).Value = 2
$ws.Cells.Item(103, 8).Value = 83

# Row 104: now "Guinea Ecuatorial"
$ws.Cells.Item(104, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(104, 2).Value = 3071
$ws.Cells.Item(104, 4).Value = 842
$ws.Cells.Item(104, 5).Value = 2178
$ws.Cells.Item(104, 8).Value = 51

# Row 112: refreshed totals
$ws.Cells.Item(112, 2).Value = 2084
$ws.Cells.Item(112, 3).Value = 3
$ws.Cells.Item(112, 4).Value = 1967
$ws.Cells.Item(112, 5).Value = 106

# Row 118: now "Malaui"
$ws.Cells.Item(118, 1).Value = "Malaui"
$ws.Cells.Item(118, 2).Value = 1864
$ws.Cells.Item(118, 3).Value = 46
$ws.Cells.Item(118, 4).Value = 345
$ws.Cells.Item(118, 5).Value = 1495
$ws.Cells.Item(118, 7).Value = 5
$ws.Cells.Item(118, 8).Value = 24

# Row 119: now "Lituania"
$ws.Cells.Item(119, 1).Value = "Lituania"
$ws.Cells.Item(119, 2).Value = 1854
$ws.Cells.Item(119, 3).Value = 10
$ws.Cells.Item(119, 4).Value = 1552
$ws.Cells.Item(119, 5).Value = 223
$ws.Cells.Item(119, 8).Value = 79

# Row 209: now "Islas Malvinas"
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"

# Row 210: now "Groenlandia"
$ws.Cells.Item(210, 1).Value = "Groenlandia"
